# Update the YouTube hyperlink caption text that appears on the "button
# toggle" video slides. The visible text is changed from the old video
# URL (N_KxVoq2HyM, index=5) to the new one (cRvqC9hPbLI, index=6); the
# hyperlink relationship itself (rId2 / rId3) is left untouched, matching
# the author's edit.

$oldUrl = "https://www.youtube.com/watch?v=N_KxVoq2HyM&list=PLC3y8-rFHvwilEuCqFGTL5Gt5U6deIrsU&index=5"
$newUrl = "https://www.youtube.com/watch?v=cRvqC9hPbLI&list=PLC3y8-rFHvwilEuCqFGTL5Gt5U6deIrsU&index=6"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $textRange = $shape.TextFrame.TextRange
            if ($textRange.Text -eq $oldUrl) {
                $textRange.Text = $newUrl
            }
        }
    }
}
